$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C (Förändrad) for rows 2..500 from 45181 to 45182
$ws.Range("C2:C500").Value2 = 45182

# 2. Ensure row 500 has an explicit row height (15), matching the other data rows
$ws.Rows.Item(500).RowHeight = 15

# 3. Add new row 501 with data
$ws.Range("A501").Value2 = "A 42524-2023"

$ws.Range("B501").Value2 = 45181
$ws.Range("B501").NumberFormat = "YYYY-MM-DD"

$ws.Range("C501").Value2 = 45182
$ws.Range("C501").NumberFormat = "YYYY-MM-DD"

$ws.Range("D501").Value2 = "HALLANDS LÄN"
$ws.Range("E501").Value2 = "FALKENBERG"

$ws.Range("G501").Value2 = 1.1
$ws.Range("H501").Value2 = 0
$ws.Range("I501").Value2 = 0
$ws.Range("J501").Value2 = 0
$ws.Range("K501").Value2 = 0
$ws.Range("L501").Value2 = 0
$ws.Range("M501").Value2 = 0
$ws.Range("N501").Value2 = 0
$ws.Range("O501").Value2 = 0
$ws.Range("P501").Value2 = 0
$ws.Range("Q501").Value2 = 0

# R501 stays empty but carries the wrap-text style used throughout column R
$ws.Range("R501").WrapText = $true
